$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-05 Thursday", "2024-12-06 Friday"),
    @("919×4=", "355×8="),
    @("375×3=", "770×5="),
    @("222×9=", "453×4="),
    @("794×6=", "586×7="),
    @("368×2=", "922×4="),
    @("230×3=", "153×4="),
    @("601×6=", "763×8="),
    @("641×6=", "272×3="),
    @("561×6=", "547×8="),
    @("757×8=", "152×6="),
    @("463×7=", "881×2="),
    @("255×5=", "210×8="),
    @("603×2=", "774×3="),
    @("797×6=", "190×5="),
    @("689×6=", "614×8="),
    @("415×6=", "913×9="),
    @("257×3=", "232×8="),
    @("279×2=", "663×9="),
    @("889×2=", "722×2="),
    @("137×4=", "604×4="),
    @("943×8=", "649×8="),
    @("851×8=", "899×7="),
    @("568×3=", "902×5="),
    @("886×5=", "320×4="),
    @("914×5=", "450×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
